# edit.ps1 - applies the "student response + parent letter" edit:
#   1. "{name}" -> "{firstname} {lastname}" in the opening paragraph, with the
#      same run-splitting Word itself produces when a user edits piecemeal,
#      and relocates the "_GoBack" bookmark into that paragraph (right after
#      ", taking up ").
#   2. Removes the signature line ("___________________________") and the
#      "Parent's/Guardian's Signature" label paragraph entirely.
#   3. Merges the "Parent's/Guardian's" + " Name" runs into a single run.

$d = $word.ActiveDocument

function Find-Range($searchText) {
    $r = $d.Range(0, $d.Content.End)
    $ok = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw ("Text not found: " + $searchText)
    }
    return $r
}

function Force-Split($range) {
    # Forces a run boundary at the edges of $range without changing its
    # formatting (toggle a character property on then back off).
    $range.Bold = $true
    $range.Bold = $false
}

# ---------------------------------------------------------------------------
# 1) {name} -> {firstname} {lastname}
# ---------------------------------------------------------------------------

# 1a. Replace the inner "name" token with "firstname" (keeps braces intact).
$nameTok = Find-Range("{name}")
$innerName = $d.Range($nameTok.Start + 1, $nameTok.End - 1)
$innerName.Text = "firstname"

# 1b. Insert " {lastname}" right after the closing brace of "{firstname}".
$firstnameTok = Find-Range("{firstname}")
$afterBrace = $d.Range($firstnameTok.End, $firstnameTok.End)
$afterBrace.InsertAfter(" {lastname}")

# 1c. Re-split the run boundaries so the paragraph matches Word's own
#     piecemeal-edit run layout:
#       "...Mr./Ms. {" | "firstname" | "}" | " {lastname}" | ", taking up " | "{course}..."
$anchor = Find-Range("Mr./Ms. {firstname} {lastname}, taking up {course}")
$pos = $anchor.Start
$openBrace      = $pos + 8                 # length of "Mr./Ms. "
$firstnameStart = $openBrace + 1
$firstnameEnd   = $firstnameStart + 9       # length of "firstname"
$closeBrace     = $firstnameEnd
$lastBlockStart = $closeBrace + 1
$lastBlockEnd   = $lastBlockStart + 11      # length of " {lastname}"
$takingUpStart  = $lastBlockEnd
$takingUpEnd    = $takingUpStart + 12       # length of ", taking up "

Force-Split ($d.Range($openBrace, $firstnameStart))        # "{"
Force-Split ($d.Range($firstnameStart, $firstnameEnd))      # "firstname"
Force-Split ($d.Range($firstnameEnd, $lastBlockStart))      # "}"
Force-Split ($d.Range($lastBlockStart, $lastBlockEnd))      # " {lastname}"
Force-Split ($d.Range($takingUpStart, $takingUpEnd))        # ", taking up "

# 1d. Relocate the "_GoBack" bookmark to sit right after ", taking up "
#     (adding a bookmark with an existing name moves it).
$bookmarkSpot = $d.Range($takingUpEnd, $takingUpEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

# ---------------------------------------------------------------------------
# 2) Remove the signature-line paragraphs.
# ---------------------------------------------------------------------------

$sigLine = Find-Range("___________________________")
$sigLine.Paragraphs.Item(1).Range.Delete()

$sigLabel = Find-Range("Parent’s/Guardian’s Signature")
$sigLabel.Paragraphs.Item(1).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Merge "Parent's/Guardian's" + " Name" into a single run.
# ---------------------------------------------------------------------------

$nameLabel = Find-Range("Parent’s/Guardian’s Name")
$nameLabel.Text = "Parent’s/Guardian’s Name"
